$wb = $excel.ActiveWorkbook

# --- survey sheet: drop "::language" suffix from guidance_hint / media::* columns ---
$survey = $wb.Worksheets.Item("survey")
$survey.Range("E1").Value = "guidance_hint"
$survey.Range("X1").Value = "media::image"
$survey.Range("Y1").Value = "media::video"
$survey.Range("Z1").Value = "media::audio"

# --- choices sheet: same trim for the media::* columns ---
$choices = $wb.Worksheets.Item("choices")
$choices.Range("D1").Value = "media::image"
$choices.Range("E1").Value = "media::video"
$choices.Range("F1").Value = "media::audio"
